$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B column (CRS Number) corrections ---------------------------------
$ws.Range("B2").Value = 59179638
$ws.Range("B5").Value = 80862497

# --- Add the new "Description" column (E) ------------------------------
# Column D ("Booking.com Price") previously held the booking status text
# ("Canceled" / "Checked in, but different date"). That status now moves
# into a new column E, and column D is populated with the actual price.
# Re-set C as well so every row is rewritten consistently.

# Row 2 - Fedrique PIERRE
$ws.Range("C2").Value = "Fedrique PIERRE"
$ws.Range("D2").Value = "242.97 USD"
$ws.Range("E2").Value = "Canceled"

# Row 3 - Yaning Qiao
$ws.Range("C3").Value = "Yaning Qiao"
$ws.Range("D3").Value = "118.15 USD"
$ws.Range("E3").Value = "Checked in, but different date"

# Row 4 - dongxu liu
$ws.Range("C4").Value = "dongxu liu"
$ws.Range("D4").Value = "95.45 USD"
$ws.Range("E4").Value = "Checked in, but different date"

# Row 5 - Terrance Graham
$ws.Range("C5").Value = "Terrance Graham"
$ws.Range("D5").Value = "119 USD"
$ws.Range("E5").Value = "Canceled"

# Row 6 - dongxu liu
$ws.Range("C6").Value = "dongxu liu"
$ws.Range("D6").Value = "98.77 USD"
$ws.Range("E6").Value = "Checked in, but different date"
